$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2.2.8 Co-Curricular and Extra C")

# Row 8 gets its first real data entry.
$ws.Range("B8").Value2 = 1

$ws.Range("C8").Value2 = 'How to Write on Wikipedia in "Marathi Language"'

# Header cell E7: "Conducting authority" -> "Resource Person/Conducting authority"
$ws.Range("E7").Value2 = "Resource Person/Conducting authority"

$ws.Range("D8").Value2 = 43843
$ws.Range("D8").NumberFormat = "mm-dd-yy"

$ws.Range("F8").Value2 = 61

$ws.Range("G8").Value2 = "Intitute Level Activity"

$ws.Range("E8").Value2 = "Mrs. A S Patil, Lecturer in Computer Engg, G.R.W.P Latur"

$ws.Range("C8").WrapText = $true
$ws.Range("E8").WrapText = $true

# Row height grows to fit the wrapped text.
$ws.Range("B8:G8").RowHeight = 25.5

# Column C widens to fit the new, longer content.
$ws.Columns("C:C").ColumnWidth = 46.85546875

# Selection moves on to the next cell to fill in.
$ws.Range("E9").Select()
